$wb = $excel.ActiveWorkbook

# Rename "Sheet3" to "nbaCityNameTeamName"
$sheet3 = $wb.Worksheets.Item("Sheet3")
$sheet3.Name = "nbaCityNameTeamName"

# myTrips: deselect A12, select A2:A17 (active cell A2), remove tabSelected
$myTrips = $wb.Worksheets.Item("myTrips")
$myTrips.Range("A2:A17").Select()

# nbaCityNameTeamName: select D24 (active cell), and make it the active sheet/tab
$sheet3.Activate()
$sheet3.Range("D24").Select()
